$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1500
$ws.Range("I12").Value = 1500
$ws.Range("K12").Value = 1500
$ws.Range("M12").Value = -1330

$ws.Range("H137").Value = 3403
$ws.Range("I137").Value = 2779.5715
$ws.Range("K137").Value = 8338.7145
$ws.Range("M137").Value = -5788.7145

$ws.Range("H138").Value = 5424.7144
$ws.Range("J138").Value = 6591.3887
$ws.Range("L138").Value = 19774.1661
$ws.Range("N138").Value = -30054.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 31999.5
$ws.Range("J16").Value = 31999.5
$ws.Range("L16").Value = 31999.5
$ws.Range("N16").Value = -32573.5

$ws.Range("H17").Value = 21332.334
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H19").Value = 30004
$ws.Range("I19").Value = 30004
$ws.Range("K19").Value = 30004
$ws.Range("M19").Value = -29775

$ws.Range("H27").Value = 28500
$ws.Range("J27").Value = 28500
$ws.Range("L27").Value = 28500
$ws.Range("N27").Value = -28868

$ws.Range("H30").Value = 26500
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 26500
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").Value = 26500
$ws.Range("N30").Value = -26800

$ws.Range("H32").Value = 8399.206
$ws.Range("I32").Value = 3083.1936
$ws.Range("K32").Value = 3083.1936
$ws.Range("M32").Value = -2796.1936

$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("K33").Value = 10000
$ws.Range("M33").Value = -9671

$ws.Range("H36").Value = 1205.2
$ws.Range("I36").Value = 1205.2
$ws.Range("K36").Value = 1205.2
$ws.Range("M36").Value = -859.2

$ws.Range("H110").Value = 1050.2646
$ws.Range("I110").Value = 1050.2646
$ws.Range("K110").Value = 1050.2646
$ws.Range("M110").Value = 994.7354

$ws.Range("H122").Value = 2702.4736
$ws.Range("I122").Value = 2643.7354
$ws.Range("K122").Value = 7931.206200000001
$ws.Range("M122").Value = -5481.206200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 9998.5
$ws.Range("I16").Value = 9998
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 9998
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -9828
$ws.Range("N16").Value = -10339

$ws.Range("H17").Value = 4474.5
$ws.Range("J17").Value = 4474.5
$ws.Range("L17").Value = 4474.5
$ws.Range("N17").Value = -4818.5

$ws.Range("H18").Value = 5011
$ws.Range("J18").Value = 5011
$ws.Range("L18").Value = 5011
$ws.Range("N18").Value = -6069

$ws.Range("H22").Value = 416.26086
$ws.Range("I22").Value = 404.8125
$ws.Range("K22").Value = 404.8125
$ws.Range("M22").Value = -231.8125

$ws.Range("H107").Value = 957.25
$ws.Range("I107").Value = 979.7143
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 979.7143
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 940.2857
$ws.Range("N107").Value = -4640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1101
$ws.Range("I22").Value = 167.57143
$ws.Range("J22").Value = 1917.75
$ws.Range("K22").Value = 167.57143
$ws.Range("L22").Value = 1917.75
$ws.Range("M22").Value = 182.42857
$ws.Range("N22").Value = -2617.75

$ws.Range("H99").Value = 2157.4583
$ws.Range("I99").Value = 1984.8857
$ws.Range("J99").Value = 2622.077
$ws.Range("K99").Value = 1984.8857
$ws.Range("L99").Value = 2622.077
$ws.Range("M99").Value = -486.8857
$ws.Range("N99").Value = -5618.077

$ws.Range("H122").Value = 1730.6923
$ws.Range("J122").Value = 2999.5
$ws.Range("L122").Value = 8998.5
$ws.Range("N122").Value = -13898.5

$ws.Range("H126").Value = 2157.4583
$ws.Range("I126").Value = 1984.8857
$ws.Range("J126").Value = 2622.077
$ws.Range("K126").Value = 5954.6571
$ws.Range("L126").Value = 7866.231000000001
$ws.Range("M126").Value = -3484.6571
$ws.Range("N126").Value = -12806.231

$ws.Range("H132").Value = 4820.8335
$ws.Range("I132").Value = 2598.3845
$ws.Range("K132").Value = 7795.1535
$ws.Range("M132").Value = -5265.1535

$ws.Range("H134").Value = 4460.636
$ws.Range("I134").Value = 2838.5715
$ws.Range("J134").Value = 7299.25
$ws.Range("K134").Value = 8515.7145
$ws.Range("L134").Value = 21897.75
$ws.Range("M134").Value = -5980.7145
$ws.Range("N134").Value = -26967.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 325.25
$ws.Range("I8").Value = 325.25
$ws.Range("K8").Value = 975.75
$ws.Range("M8").Value = -836.75

$ws.Range("H12").Value = 340.3889
$ws.Range("J12").Value = 436.42856
$ws.Range("L12").Value = 1309.28568
$ws.Range("N12").Value = -1655.28568

$ws.Range("H97").Value = 1428.5
$ws.Range("I97").Value = 3900
$ws.Range("J97").Value = 604.6667
$ws.Range("K97").Value = 11700
$ws.Range("L97").Value = 1814.0001
$ws.Range("M97").Value = -11204
$ws.Range("N97").Value = -2806.0001

$ws.Range("H98").Value = 843.375
$ws.Range("J98").Value = 688.8
$ws.Range("L98").Value = 2066.4
$ws.Range("N98").Value = -5062.4

$ws.Range("H131").Value = 759268.2
$ws.Range("I131").Value = 954.4167
$ws.Range("J131").Value = 1517581.9
$ws.Range("K131").Value = 2863.2501
$ws.Range("L131").Value = 4552745.699999999
$ws.Range("M131").Value = 2176.7499
$ws.Range("N131").Value = -4562825.699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 788.3333
$ws.Range("I107").Value = 557.5
$ws.Range("K107").Value = 557.5
$ws.Range("M107").Value = 1362.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83122.766
$ws.Range("I7").Value = 117010.89
$ws.Range("J7").Value = 6874.5
$ws.Range("K7").Value = 117010.89
$ws.Range("L7").Value = 6874.5
$ws.Range("M7").Value = -116898.89
$ws.Range("N7").Value = -7098.5

$ws.Range("H16").Value = 1543.2222
$ws.Range("J16").Value = 1454.8334
$ws.Range("L16").Value = 1454.8334
$ws.Range("N16").Value = -1794.8334

$ws.Range("H22").Value = 1856.5
$ws.Range("I22").Value = 532
$ws.Range("J22").Value = 2651.2
$ws.Range("K22").Value = 532
$ws.Range("L22").Value = 2651.2
$ws.Range("M22").Value = -237
$ws.Range("N22").Value = -3241.2

$ws.Range("H27").Value = 1856.5
$ws.Range("I27").Value = 532
$ws.Range("J27").Value = 2651.2
$ws.Range("K27").Value = 532
$ws.Range("L27").Value = 2651.2
$ws.Range("M27").Value = -425
$ws.Range("N27").Value = -2865.2

$ws.Range("H46").Value = 4514.615
$ws.Range("J46").Value = 4649.125
$ws.Range("L46").Value = 4649.125
$ws.Range("N46").Value = -5025.125

$ws.Range("H105").Value = 130756.25
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H122").Value = 6305.6665
$ws.Range("I122").Value = 5066.9
$ws.Range("K122").Value = 15200.7
$ws.Range("M122").Value = -12750.7

$ws.Range("H126").Value = 83122.766
$ws.Range("I126").Value = 117010.89
$ws.Range("J126").Value = 6874.5
$ws.Range("K126").Value = 351032.67
$ws.Range("L126").Value = 20623.5
$ws.Range("M126").Value = -348562.67
$ws.Range("N126").Value = -25563.5

$ws.Range("H132").Value = 6797.68
$ws.Range("I132").Value = 6375.25
$ws.Range("K132").Value = 19125.75
$ws.Range("M132").Value = -16595.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 147820.67
$ws.Range("J64").Value = 147820.67
$ws.Range("L64").Value = 147820.67
$ws.Range("N64").Value = -148316.67

$ws.Range("H67").Value = 147820.67
$ws.Range("J67").Value = 147820.67
$ws.Range("L67").Value = 147820.67
$ws.Range("N67").Value = -149536.67

$ws.Range("H136").Value = 6654.864
$ws.Range("I136").Value = 5766.647
$ws.Range("J136").Value = 9674.799999999999
$ws.Range("K136").Value = 17299.941
$ws.Range("L136").Value = 29024.4
$ws.Range("M136").Value = -14749.941
$ws.Range("N136").Value = -34124.39999999999
